$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Refresh the crypto symbol list snapshot (GitHub Actions bot commit) ---

# Every data row (2-51) moves forward one hour: Hora (column G) "18" -> "19".
# Set Text format first so the value keeps being stored as a string, matching the source sheet.
$gRange = $ws.Range("G2:G51")
$gRange.NumberFormat = "@"
$gRange.Value = "19"

# Row-specific Price (column D) refreshes.
$dRows = @(2, 3, 4, 5, 6, 7, 8, 9, 10, 11, 12, 13, 14, 15, 16, 17, 18, 20, 22, 23, 24, 25, 26, 27, 28, 40, 41, 43, 44, 45, 47)
$dVals = @("287.16", "29.20", "5.205", "0.06965", "7.435", "3.556", "1.405", "0.9007", "0.1612", "0.07515", "0.07767", "0.02909", "0.09017", "0.001587", "0.0006494", "0.006096", "3.469", "0.3246", "4.037", "0.1600", "0.04523", "0.001211", "0.004142", "0.0001170", "0.0001644", "0.04388", "0.006918", "0.002069", "0.01158", "0.00005859", "0.01300")
for ($i = 0; $i -lt $dRows.Count; $i++) {
    $cell = $ws.Range("D" + $dRows[$i])
    $cell.NumberFormat = "@"
    $cell.Value = $dVals[$i]
}

# Row-specific Volume(1h) (column E) refreshes.
$eRows = @(2, 3, 4, 5, 6, 7, 8, 9, 10, 11, 12, 13, 14, 15, 16, 17, 18, 19, 20, 21, 22, 23, 24, 25, 26, 27, 28, 40, 41, 42, 43, 44, 45, 46, 47)
$eVals = @("0.76%", "2.32%", "2.17%", "4.54%", "1.42%", "5.03%", "3.20%", "-3.80%", "2.77%", "15.83%", "1.47%", "0.89%", "0.59%", "0.02%", "1.18%", "-0.90%", "-0.23%", "0.13%", "1.44%", "2.26%", "-0.36%", "5.08%", "1.21%", "2.72%", "-7.35%", "-6.08%", "1.95%", "5.02%", "2.73%", "-0.16%", "2.83%", "0.30%", "3.51%", "-1.86%", "-0.21%")
for ($i = 0; $i -lt $eRows.Count; $i++) {
    $cell = $ws.Range("E" + $eRows[$i])
    $cell.NumberFormat = "@"
    $cell.Value = $eVals[$i]
}
